# Edit script: add Google-Sheets style attendance export columns for 2025-11-22,
# refresh running totals for existing students, and append two newly joined students.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header columns for the 2025-11-22 session (X1:AD1), reusing the
#    bold/centered/bordered header look that is already applied to A1:W1
#    (style index 1). We copy formats (not values) from W1 across the new
#    header cells, both before *and* after filling in their text, so that
#    the trailing "2025-11-22" header - which looks like a date - keeps the
#    literal text value instead of being auto-converted to a date serial.
# ---------------------------------------------------------------------------
$headerSrc = $ws.Range("W1")
$newHeaderRange = $ws.Range("X1:AD1")

$headerSrc.Copy()
$newHeaderRange.PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("X1").Value = "2025-11-22_x"
$ws.Range("Y1").Value = "2025-11-22_y"
$ws.Range("Z1").Value = "2025-11-22_x.1"
$ws.Range("AA1").Value = "2025-11-22_y.1"
$ws.Range("AB1").Value = "2025-11-22_x.2"
$ws.Range("AC1").Value = "2025-11-22_y.2"
# Leading apostrophe forces this to stay plain text rather than becoming a date.
$ws.Range("AD1").Value = "'2025-11-22"

# Re-copy the header format on top so the quote-prefix text entry above ends
# up sharing the exact same style as the rest of the header row.
$headerSrc.Copy()
$newHeaderRange.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Existing students (rows 2-8): the session total grew from 17 to 22, so
#    refresh the "Total" column and mark the new 2025-11-22 slots absent.
#    Row 7 (Shubham Phad) keeps 3 "Present" days, so her Attendance % drops
#    from 17.6 to 13.6 once the new total is applied.
# ---------------------------------------------------------------------------
$newColumns = @("X", "Y", "Z", "AA", "AB", "AC", "AD")

for ($row = 2; $row -le 8; $row++) {
    $ws.Range("V$row").Value = 22

    foreach ($col in $newColumns) {
        $ws.Range("$col$row").Value = "❌"
    }
}

$ws.Range("W7").Value = 13.6

# ---------------------------------------------------------------------------
# 3. Two newly onboarded students (rows 9-10). They have no attendance marks
#    for the earlier 2025-11-07 .. 2025-11-21 sessions (columns D:T stay
#    blank) and only 4 sessions worth of "Total" so far.
# ---------------------------------------------------------------------------
$newStudents = @(
    @{ Row = 9;  Roll = "EC4226"; Name = "Abhishek Pathak"; Email = "abhipathak2513@gmail.com" },
    @{ Row = 10; Roll = "EC4237"; Name = "Shubham Pitekar"; Email = "shubhampitekar2323@gmail.com" }
)

foreach ($student in $newStudents) {
    $row = $student.Row

    $ws.Range("A$row").Value = $student.Roll
    $ws.Range("B$row").Value = $student.Name
    $ws.Range("C$row").Value = $student.Email

    $ws.Range("U$row").Value = 0
    $ws.Range("V$row").Value = 4
    $ws.Range("W$row").Value = 0

    # Row 9 has no mark yet for the very first new-session column (X);
    # every other new column (and all of the legacy D:T columns) is absent.
    foreach ($col in $newColumns) {
        if (-not ($row -eq 9 -and $col -eq "X")) {
            $ws.Range("$col$row").Value = "❌"
        }
    }
}

